{"js": "const newValues = [\n  \"8-7=1\",\n  \"22+71=93\",\n  \"96-82=14\",\n  \"57+40=97\",\n  \"42-3=39\",\n  \"52-3=49\",\n  \"14+70=84\",\n  \"71+16=87\",\n  \"77+0=77\",\n  \"8+60=68\",\n  \"62-27=35\",\n  \"67-43=24\",\n  \"65-30=35\",\n  \"16+68=84\",\n  \"16+16=32\",\n  \"77-24=53\",\n  \"83-67=16\",\n  \"1+29=30\",\n  \"32+64=96\",\n  \"72-19=53\",\n  \"2+69=71\",\n  \"51-50=1\",\n  \"4+12=16\",\n  \"93-63=30\",\n  \"19+3=22\",\n  \"18+57=75\",\n  \"72-46=26\",\n  \"64-8=56\",\n  \"77+20=97\",\n  \"91+1=92\",\n  \"31-5=26\",\n  \"33+15=48\",\n  \"70-39=31\",\n  \"12+74=86\",\n  \"54-11=43\",\n  \"29-4=25\",\n  \"13+21=34\",\n  \"85-60=25\",\n  \"95-44=51\",\n  \"93-3=90\",\n  \"40+24=64\",\n  \"51-3=48\",\n  \"0+30=30\",\n  \"13-2=11\",\n  \"51+11=62\",\n  \"98-62=36\",\n  \"21-3=18\",\n  \"67+32=99\",\n  \"8+68=76\",\n  \"86-65=21\",\n  \"34+37=71\",\n  \"74-53=21\",\n  \"15+22=37\",\n  \"1+39=40\",\n  \"38-9=29\",\n  \"89-1=88\",\n  \"43+17=60\",\n  \"52-17=35\",\n  \"54-45=9\",\n  \"13+34=47\",\n  \"84-43=41\",\n  \"41+31=72\",\n  \"78-73=5\",\n  \"2+87=89\",\n  \"4+27=31\",\n  \"59-6=53\",\n  \"21+58=79\",\n  \"34+15=49\",\n  \"57-20=37\",\n  \"40-6=34\",\n  \"9+80=89\",\n  \"93-9=84\",\n  \"16+62=78\",\n  \"3+91=94\",\n  \"16+60=76\",\n  \"73-17=56\",\n  \"27+44=71\",\n  \"4+6=10\",\n  \"67-44=23\",\n  \"13+44=57\",\n  \"98-39=59\",\n  \"56-29=27\",\n  \"70+6=76\",\n  \"27+59=86\",\n  \"6+43=49\",\n  \"28+65=93\",\n  \"73-68=5\",\n  \"41+49=90\",\n  \"28+67=95\",\n  \"61-0=61\",\n  \"47+38=85\",\n  \"58-24=34\",\n  \"77-51=26\",\n  \"47+9=56\",\n  \"66+26=92\",\n  \"31-13=18\",\n  \"65-39=26\",\n  \"24+5=29\",\n  \"39-15=24\",\n  \"87-85=2\"\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The document has a single table holding all the practice-problem cells.\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\n// Replacements are positional (row-major): some original expressions repeat\n// (e.g. \"92-74=18\" appears twice) but map to different replacements at\n// different positions, so we must walk cell-by-cell rather than doing a\n// global text search/replace.\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "$newValues = @(\n    \"8-7=1\",\n    \"22+71=93\",\n    \"96-82=14\",\n    \"57+40=97\",\n    \"42-3=39\",\n    \"52-3=49\",\n    \"14+70=84\",\n    \"71+16=87\",\n    \"77+0=77\",\n    \"8+60=68\",\n    \"62-27=35\",\n    \"67-43=24\",\n    \"65-30=35\",\n    \"16+68=84\",\n    \"16+16=32\",\n    \"77-24=53\",\n    \"83-67=16\",\n    \"1+29=30\",\n    \"32+64=96\",\n    \"72-19=53\",\n    \"2+69=71\",\n    \"51-50=1\",\n    \"4+12=16\",\n    \"93-63=30\",\n    \"19+3=22\",\n    \"18+57=75\",\n    \"72-46=26\",\n    \"64-8=56\",\n    \"77+20=97\",\n    \"91+1=92\",\n    \"31-5=26\",\n    \"33+15=48\",\n    \"70-39=31\",\n    \"12+74=86\",\n    \"54-11=43\",\n    \"29-4=25\",\n    \"13+21=34\",\n    \"85-60=25\",\n    \"95-44=51\",\n    \"93-3=90\",\n    \"40+24=64\",\n    \"51-3=48\",\n    \"0+30=30\",\n    \"13-2=11\",\n    \"51+11=62\",\n    \"98-62=36\",\n    \"21-3=18\",\n    \"67+32=99\",\n    \"8+68=76\",\n    \"86-65=21\",\n    \"34+37=71\",\n    \"74-53=21\",\n    \"15+22=37\",\n    \"1+39=40\",\n    \"38-9=29\",\n    \"89-1=88\",\n    \"43+17=60\",\n    \"52-17=35\",\n    \"54-45=9\",\n    \"13+34=47\",\n    \"84-43=41\",\n    \"41+31=72\",\n    \"78-73=5\",\n    \"2+87=89\",\n    \"4+27=31\",\n    \"59-6=53\",\n    \"21+58=79\",\n    \"34+15=49\",\n    \"57-20=37\",\n    \"40-6=34\",\n    \"9+80=89\",\n    \"93-9=84\",\n    \"16+62=78\",\n    \"3+91=94\",\n    \"16+60=76\",\n    \"73-17=56\",\n    \"27+44=71\",\n    \"4+6=10\",\n    \"67-44=23\",\n    \"13+44=57\",\n    \"98-39=59\",\n    \"56-29=27\",\n    \"70+6=76\",\n    \"27+59=86\",\n    \"6+43=49\",\n    \"28+65=93\",\n    \"73-68=5\",\n    \"41+49=90\",\n    \"28+67=95\",\n    \"61-0=61\",\n    \"47+38=85\",\n    \"58-24=34\",\n    \"77-51=26\",\n    \"47+9=56\",\n    \"66+26=92\",\n    \"31-13=18\",\n    \"65-39=26\",\n    \"24+5=29\",\n    \"39-15=24\",\n    \"87-85=2\"\n)\n\n# Replacements are positional (row-major): some original expressions repeat\n# (e.g. \"92-74=18\" appears twice) but map to different replacements at\n# different positions, so we must walk cell-by-cell in table order rather\n# than doing a global Find/Replace on the text.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $row = $t.Rows.Item($r)\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $cell = $row.Cells.Item($c)\n        $rng = $cell.Range\n        # Exclude the trailing cell-mark character so we don't clobber it\n        $rng.End = $rng.End - 1\n        $rng.Text = $newValues[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
